$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day 4 header date gets filled in (column E, row 4) ---
$ws.Range("E4").Value = "第四天`n日期:2025-11-25"

# Day 3 (D column) "lights-out" time correction, plus Day 4 value
$ws.Range("D8").Value = "24：00"
$ws.Range("E8").Value = "23：00"

# --- Day 4 column (E) sleep-diary entries ---
$ws.Range("E5").Value = "8：33"
$ws.Range("E6").Value = "8：40"
$ws.Range("E7").Value = "22：20"

$ws.Range("E9").Value = 30
$ws.Range("E10").Value = 4
$ws.Range("E11").Value = 30
$ws.Range("E12").Value = 450
$ws.Range("E13").Value = "无"
$ws.Range("E14").Value = 40
$ws.Range("E15").Value = 3
$ws.Range("E16").Value = 4
$ws.Range("E17").Value = 2
$ws.Range("E18").Value = "无"

# --- Restore the view / selection state recorded at save time ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E18").Select()
